$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Japhd932"
$ws.Range("B2").Value = 23082104
$ws.Range("C2").Value = "gihbmrn38"
$ws.Range("D2").Value = "pMD&#32j"
$ws.Range("F2").Value = "RymAsYwy"
$ws.Range("G2").Value = "jEHF"
